$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrigiendo permisos del Regente: quitar "formulas" de la lista de permisos
$ws.Range("B3").Value = "CRUD(medicamentos,sucursales)"

# Reflejar la posicion final de seleccion de la celda tras la edicion
$ws.Range("B18").Select()
